$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updated values for columns D (Price), E (Volume(1h)), G (Hora).
# Only columns that actually changed for a given row are listed.
$updates = @(
    @{ Row = 2; D = "309.21"; E = "-3.85%"; G = "7" },
    @{ Row = 3; D = "49.73"; E = "1.47%"; G = "7" },
    @{ Row = 4; D = "5.180"; E = "-2.09%"; G = "7" },
    @{ Row = 5; D = "0.07753"; E = "-4.34%"; G = "7" },
    @{ Row = 6; D = "4.527"; E = "-2.00%"; G = "7" },
    @{ Row = 7; D = "1.367"; E = "12.53%"; G = "7" },
    @{ Row = 8; D = "1.554"; E = "-6.64%"; G = "7" },
    @{ Row = 9; D = "0.1238"; E = "-6.19%"; G = "7" },
    @{ Row = 10; D = "0.1959"; E = "0.10%"; G = "7" },
    @{ Row = 11; D = "0.04714"; E = "4.55%"; G = "7" },
    @{ Row = 12; D = "0.09357"; E = "-1.44%"; G = "7" },
    @{ Row = 13; D = "0.1045"; E = "-0.14%"; G = "7" },
    @{ Row = 14; D = "0.001259"; E = "-5.14%"; G = "7" },
    @{ Row = 15; D = "0.04170"; E = "-3.46%"; G = "7" },
    @{ Row = 16; D = "0.005807"; E = "-1.29%"; G = "7" },
    @{ Row = 17; E = "2,016.22%"; G = "7" },
    @{ Row = 18; E = "-0.82%"; G = "7" },
    @{ Row = 19; D = "2.237"; E = "-8.24%"; G = "7" },
    @{ Row = 20; E = "2.88%"; G = "7" },
    @{ Row = 21; D = "7.924"; E = "-4.02%"; G = "7" },
    @{ Row = 22; D = "0.1341"; E = "-5.08%"; G = "7" },
    @{ Row = 23; D = "0.3040"; E = "4.12%"; G = "7" },
    @{ Row = 24; D = "0.001271"; E = "-2.93%"; G = "7" },
    @{ Row = 25; D = "0.004015"; E = "-5.55%"; G = "7" },
    @{ Row = 26; D = "0.0001350"; E = "-0.37%"; G = "7" },
    @{ Row = 27; G = "7" },
    @{ Row = 28; G = "7" },
    @{ Row = 29; G = "7" },
    @{ Row = 30; G = "7" },
    @{ Row = 31; G = "7" },
    @{ Row = 32; G = "7" },
    @{ Row = 33; G = "7" },
    @{ Row = 34; G = "7" },
    @{ Row = 35; G = "7" },
    @{ Row = 36; G = "7" },
    @{ Row = 37; G = "7" },
    @{ Row = 38; D = "0.02601"; E = "-4.85%"; G = "7" },
    @{ Row = 39; D = "0.05864"; E = "4.77%"; G = "7" },
    @{ Row = 40; D = "0.01072"; E = "69.84%"; G = "7" },
    @{ Row = 41; D = "0.007935"; E = "2.49%"; G = "7" },
    @{ Row = 42; D = "0.1420"; E = "-1.63%"; G = "7" },
    @{ Row = 43; D = "0.008456"; E = "9.79%"; G = "7" },
    @{ Row = 44; D = "0.007694"; E = "-5.03%"; G = "7" },
    @{ Row = 45; D = "0.3390"; E = "6.13%"; G = "7" },
    @{ Row = 46; D = "0.00007038"; E = "0.75%"; G = "7" },
    @{ Row = 47; E = "0.01%"; G = "7" },
    @{ Row = 48; D = "0.04870"; E = "-20.59%"; G = "7" },
    @{ Row = 49; D = "0.002625"; E = "-34.53%"; G = "7" },
    @{ Row = 50; E = "0.01%"; G = "7" },
    @{ Row = 51; D = "0.0002004"; E = "0.01%"; G = "7" }
)

# Use a scratch cell formatted as Text so values such as "7" or "-3.85%"
# round-trip as literal strings (matching the inline-string cells already
# used throughout this sheet) instead of being coerced to numbers/percentages.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

foreach ($u in $updates) {
    foreach ($col in @("D", "E", "G")) {
        if ($u.ContainsKey($col)) {
            $scratch.Value = $u[$col]
            $scratch.Copy()
            $ws.Range($col + $u.Row).PasteSpecial(-4163)
        }
    }
}

# Remove the scratch cell so it leaves no trace in the used range.
$scratch.Clear()